# Apply changes described by the commit:
# "added isdriveropen method; corrected ctdc tc to run; added queries in all ctdc tc xls"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the Neo4j query text into A2 (new shared string, wraps text via existing style)
$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Serous endometrial adenocarcinoma'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

$ws.Range("A2").Value = $query

# Adjust row height for row 2 to fit the wrapped query text
$ws.Rows.Item(2).RowHeight = 87

# Update selection / view to B2 (no longer frozen at B1 topLeftCell)
$ws.Range("B2").Select() | Out-Null

$wb.Save() | Out-Null
